# Removed Test Case Inter-Dependency
#
# - Make "ProductLoanInput" the active/selected sheet (was "ProductLoanOutput"),
#   and reset its selection back to cell B1 (was B9).
# - Update the product name text (shared between the input and output sheets)
#   from "...VAR-INST-UPFRONT" to "...VAR-INST-UP1st".
# - Change the "shortname" cell on the input sheet from the numeric 2580
#   to the text value "258d".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2580-MS-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UP1st"

# Update the product name on both sheets (they share the same text) so the
# underlying shared string is edited in place rather than duplicated.
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# shortname: was the number 2580, now the text "258d"
$ws1.Range("B2").Value = "258d"

# Re-select B1 on the input sheet and make that sheet the active tab again,
# undoing the previous test run's leftover selection/tab state.
[void]$ws1.Range("B1").Select()
[void]$ws1.Activate()
